# Add variable in excel keyword
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row: rename "Object" -> "Object Name", "value" -> "Value",
#    and add a new "Using Variable" header in column F.
# ---------------------------------------------------------------------------
$ws.Range("C1").Value = "Object Name"
$ws.Range("E1").Value = "Value"

$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Using Variable"

# ---------------------------------------------------------------------------
# 2. "Reset Login In Application" block: update the username/password
#    values, then insert a new SETVARIABLE step before the CLICK step,
#    followed by a blank spacer row.
# ---------------------------------------------------------------------------
$ws.Range("E4").Value = "mngr73146"
$ws.Range("E5").Value = "dAsazYt"

# insert the new SETVARIABLE row just above the CLICK row (currently row 6)
$ws.Rows.Item(6).Insert()
$ws.Range("A6:E6").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A6:E6").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("B6").Value = "SETVARIABLE"
$ws.Range("C6").Value = "accountID"
$ws.Range("D6").Value = "xpath"

# insert a blank spacer row after the CLICK row (currently row 7, after insert above)
$ws.Rows.Item(8).Insert()
$ws.Range("A7:E7").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8:E8").Value = ""
$ws.Rows.Item(8).RowHeight = 15.75

# ---------------------------------------------------------------------------
# 3. "Login In Application" block (now shifted down to rows 9-13): the
#    SETTEXT username step now references the new "accountID" variable
#    instead of a literal value.
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = ""
$ws.Range("E11").Copy()
$ws.Range("F11").PasteSpecial(-4122)
$ws.Range("F11").Value = "accountID"

# ---------------------------------------------------------------------------
# 4. Column F width + page setup + selection, matching the edited layout.
# ---------------------------------------------------------------------------
$ws.Columns.Item(6).ColumnWidth = 13

$ws.PageSetup.Orientation = 1

$ws.Range("F12").Select()
